$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 2.57
$ws.Range("H4").Value = 2.7
$ws.Range("K4").Value = 1.77
$ws.Range("O4").Value = 1.73
$ws.Range("P4").Value = 2
$ws.Range("V4").Value = 1.5
$ws.Range("AC4").Value = 4.75

# Row 5 updates
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 6.5
$ws.Range("J5").Value = 2.05
$ws.Range("K5").Value = 2.37
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.83
$ws.Range("W5").Value = 7
$ws.Range("AG5").Value = 17
$ws.Range("AH5").Value = 34
$ws.Range("AI5").Value = 21
$ws.Range("AJ5").Value = 67
$ws.Range("AN5").Value = 7.5
$ws.Range("AR5").Value = 126
$ws.Range("AV5").Value = 7.5
$ws.Range("BA5").Value = 251
